# latest processed data after running the motilal_portfolio_change_engine
#
# Insert a new "Industry" column at column C (shifting the existing
# "Mutual Fund", "Status", "Jan_2026", "Dec_2025", "Oct_2025", "MoM", "QoQ"
# columns one position to the right, from D:J), then populate the new
# column with each holding's industry classification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:I -> D:J, leaving a blank column C to fill in.
$ws.Columns.Item(3).Insert()

# Header - use the same bold/centered style as the other header cells.
$ws.Cells.Item(1, 3).Value = "Industry"
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Row -> Industry mapping
$industries = @{
    2  = "Finance"
    3  = "Insurance"
    4  = "Construction"
    5  = "Power"
    6  = "Finance"
    7  = "Metals & Minerals Trading"
    8  = "Power"
    9  = "Finance"
    10 = "Pharmaceuticals & Biotechnology"
    11 = "Power"
    12 = "Healthcare"
    13 = "Insurance"
    14 = "Entertainment"
    15 = "Construction"
    16 = "Telecom - Services"
    17 = "Cement & Cement Products"
    18 = "Construction"
    19 = "Retailing"
    20 = "IT - Services"
    21 = "Construction"
    22 = "Pharmaceuticals & Biotechnology"
    23 = "Pharmaceuticals & Biotechnology"
    24 = "Insurance"
    25 = "Banks"
    26 = "Personal Products"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
